# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75 was previously the last row (formatted with the "last row" date
# style, i.e. the plain YYYY-MM-DD format). Since a new row is being
# appended below it, row 75 becomes a regular data row and should switch
# to the regular date/time number format used by all the other data rows.
$ws.Range("A75").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 76, keeping the "last row" date format
# (plain YYYY-MM-DD) that row 75 used to have.
$ws.Cells.Item(76, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(76, 1).Value = 45816
$ws.Cells.Item(76, 2).Value = 325
$ws.Cells.Item(76, 3).Value = 324
$ws.Cells.Item(76, 4).Value = 329
